$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Mexican Liga MX, Queretaro vs Tijuana (was Cruz Azul vs Atlas) ---
$ws.Cells.Item(2, 2).NumberFormat = "@"   # keep Date column as literal text, not an auto-converted date serial
$ws.Cells.Item(2, 1).Value = "Mexican Liga MX"
$ws.Cells.Item(2, 2).Value = "2026-01-14"
$ws.Cells.Item(2, 3).Value = "22:00:00"
$ws.Cells.Item(2, 4).Value = "Queretaro"
$ws.Cells.Item(2, 5).Value = "Tijuana"
$ws.Cells.Item(2, 6).Value = 75
$ws.Cells.Item(2, 7).Value = 290
$ws.Cells.Item(2, 8).Value = 1.05
$ws.Cells.Item(2, 9).Value = 1.06
$ws.Cells.Item(2, 10).Value = 21
$ws.Cells.Item(2, 11).Value = 26
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 1.83
$ws.Cells.Item(2, 17).Value = 2.1
$ws.Cells.Item(2, 18).Value = 1.12
$ws.Cells.Item(2, 19).Value = 8.199999999999999
$ws.Cells.Item(2, 20).Value = 1.39
$ws.Cells.Item(2, 21).Value = 1.01
$ws.Cells.Item(2, 22).Value = 18
$ws.Cells.Item(2, 23).Value = 1.01
$ws.Cells.Item(2, 24).Value = 1000
$ws.Cells.Item(2, 25).Value = 1000
$ws.Cells.Item(2, 26).Value = 1000
$ws.Cells.Item(2, 27).Value = 1000
$ws.Cells.Item(2, 28).Value = 1000
$ws.Cells.Item(2, 29).Value = 1000
$ws.Cells.Item(2, 30).Value = 990
$ws.Cells.Item(2, 31).Value = 990
$ws.Cells.Item(2, 32).Value = 1000
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 34).Value = 990
$ws.Cells.Item(2, 35).Value = 1000
$ws.Cells.Item(2, 36).Value = 1000
$ws.Cells.Item(2, 37).Value = 1000
$ws.Cells.Item(2, 38).Value = 1000
$ws.Cells.Item(2, 39).Value = 1000
$ws.Cells.Item(2, 40).Value = 1000
$ws.Cells.Item(2, 41).Value = 1000

# --- Row 3: Mexican Liga MX, CF America vs Atletico San Luis (was Queretaro vs Tijuana) ---
$ws.Cells.Item(3, 2).NumberFormat = "@"   # keep Date column as literal text, not an auto-converted date serial
$ws.Cells.Item(3, 1).Value = "Mexican Liga MX"
$ws.Cells.Item(3, 2).Value = "2026-01-14"
$ws.Cells.Item(3, 3).Value = "22:05:00"
$ws.Cells.Item(3, 4).Value = "CF America"
$ws.Cells.Item(3, 5).Value = "Atletico San Luis"
$ws.Cells.Item(3, 6).Value = 22
$ws.Cells.Item(3, 7).Value = 24
$ws.Cells.Item(3, 8).Value = 1.31
$ws.Cells.Item(3, 9).Value = 1.34
$ws.Cells.Item(3, 10).Value = 4.9
$ws.Cells.Item(3, 11).Value = 5.1
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 1.22
$ws.Cells.Item(3, 15).Value = 1.68
$ws.Cells.Item(3, 16).Value = 1.26
$ws.Cells.Item(3, 17).Value = 4.6
$ws.Cells.Item(3, 18).Value = 1.06
$ws.Cells.Item(3, 19).Value = 14
$ws.Cells.Item(3, 20).Value = 2.6
$ws.Cells.Item(3, 21).Value = 1.54
$ws.Cells.Item(3, 22).Value = 1.01
$ws.Cells.Item(3, 23).Value = 1.01
$ws.Cells.Item(3, 24).Value = 1000
$ws.Cells.Item(3, 25).Value = 2.7
$ws.Cells.Item(3, 26).Value = 5.9
$ws.Cells.Item(3, 27).Value = 1000
$ws.Cells.Item(3, 28).Value = 1000
$ws.Cells.Item(3, 29).Value = 1000
$ws.Cells.Item(3, 30).Value = 990
$ws.Cells.Item(3, 31).Value = 990
$ws.Cells.Item(3, 32).Value = 1000
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 34).Value = 990
$ws.Cells.Item(3, 35).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 37).Value = 1000
$ws.Cells.Item(3, 38).Value = 1000
$ws.Cells.Item(3, 39).Value = 1000
$ws.Cells.Item(3, 40).Value = 1000
$ws.Cells.Item(3, 41).Value = 1000

# --- Row 4 (old CF America vs Atletico San Luis row) no longer exists: remove it ---
$ws.Rows.Item(4).Delete()
